$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-apply values that already exist in the shared-string table so ---
# --- their original index is reused (order among these doesn't matter) ---
$ws.Range("B4").Value = "String"
$ws.Range("B7").Value = "String"
$ws.Range("B10").Value = "String"
$ws.Range("B13").Value = "String"
$ws.Range("B16").Value = "String"

$ws.Range("C10").Value = "array[]"
$ws.Range("C13").Value = "array^"
$ws.Range("C16").Value = "array*"

$ws.Range("B6").Value = "Datatype MyType2"
$ws.Range("B9").Value = "Datatype MyType3"
$ws.Range("B12").Value = "Datatype MyType4"
$ws.Range("B15").Value = "Datatype MyType5"

# --- Now introduce the brand-new shared strings in the exact order ---
# --- they first appear in the finished workbook ---
$ws.Range("C4").Value = "array²"
$ws.Range("C7").Value = "0array"
$ws.Range("B3").Value = "Datatype MyType1"

$ws.Range("F3").Value = "Datatype MyType1_2"
$ws.Range("F6").Value = "Datatype MyType2_2"
$ws.Range("F9").Value = "Datatype MyType3_2"
$ws.Range("F12").Value = "Datatype MyType4_2"
$ws.Range("F15").Value = "Datatype MyType5_2"

$ws.Range("J3").Value = "Datatype MyType1_3"
$ws.Range("J6").Value = "Datatype MyType2_3"
$ws.Range("J9").Value = "Datatype MyType3_3"
$ws.Range("J12").Value = "Datatype MyType4_3"
$ws.Range("J15").Value = "Datatype MyType5_3"

# --- Fill in the remaining "String" cells that mirror column B/J (reuse idx 0) ---
$ws.Range("F4").Value = "String"
$ws.Range("F7").Value = "String"
$ws.Range("F10").Value = "String"
$ws.Range("F13").Value = "String"
$ws.Range("F16").Value = "String"

$ws.Range("J4").Value = "String"
$ws.Range("J7").Value = "String"
$ws.Range("J10").Value = "String"
$ws.Range("J13").Value = "String"
$ws.Range("J16").Value = "String"

# --- Formulas for the "G" (shared-style) and "K" columns ---
$ws.Range("G4").Formula = '=C4 &" :context"'
$ws.Range("G7").Formula = '=C7 &" :context"'
$ws.Range("G10").Formula = '=C10 &" :context"'
$ws.Range("G13").Formula = '=C13 &" :context"'
$ws.Range("G16").Formula = '=C16 &" :context"'

$ws.Range("K4").Formula = '=C4 &" :context .lob"'
$ws.Range("K7").Formula = '=C7 &" :context .lob"'
$ws.Range("K10").Formula = '=C10 &" :context .lob"'
$ws.Range("K13").Formula = '=C13 &" :context .lob"'
$ws.Range("K16").Formula = '=C16 &" :context .lob"'

# --- Selection / view update ---
$ws.Range("G19").Select()
